$d = $word.ActiveDocument

# --- Change 1: merge "Front End Technologies" + " " runs into a single
#     run with text "Front End Technologies " (trailing space preserved).
$d.Content.Find.Execute("Front End Technologies ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Front End Technologies ", 2)

# --- Change 2: move "Screenshots of Running Application:" up one paragraph
#     (into the blank paragraph right after "Screenshots of Code:"), leaving
#     the paragraph that used to hold the text blank again.
$src1 = $d.Paragraphs.Item(32).Range
$dst1 = $d.Paragraphs.Item(31).Range
$dst1.FormattedText = $src1.FormattedText
$clear1 = $d.Range($src1.Start, $src1.End - 1)
$clear1.Delete()

# --- Change 3: move "URL to GitHub Repository:" up one paragraph (into the
#     blank paragraph right after "Screenshots of Running Application:" /
#     the (now) blank paragraph), leaving the paragraph that used to hold
#     the text blank again.
$src2 = $d.Paragraphs.Item(34).Range
$dst2 = $d.Paragraphs.Item(33).Range
$dst2.FormattedText = $src2.FormattedText
$clear2 = $d.Range($src2.Start, $src2.End - 1)
$clear2.Delete()

# --- Change 4: the last paragraph (now blank, still carrying the old bold
#     "URL to GitHub Repository:" formatting) becomes the new GitHub URL
#     line, styled with bCs instead of b.  Insert a brand-new paragraph with
#     the exact desired formatting right before it, then merge the old
#     (left over) blank paragraph mark away.
$lastPara = $d.Paragraphs.Item(34)
$insertAt = $lastPara.Range

$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:bCs/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:szCs w:val="24"/></w:rPr><w:t>https://github.com/mctimoth/FESD-FET-Week4</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertAt.InsertXML($xml)

$leftover = $d.Paragraphs.Item(35)
$mergeAway = $d.Range($leftover.Range.Start - 1, $leftover.Range.End)
$mergeAway.Delete()

Write-Output "done"
